# Rename the three inline picture placeholders (Pearson logo in both
# footers, BTEC logo in the first-page header) per the commit diff:
#   footer (id=3) : image1.png -> image2.png
#   footer (id=2) : image1.png -> image2.png
#   header (id=1) : image2.jpg -> image1.jpg
#
# The description (alt text) is left untouched; only the shape's Name
# (the OOXML docPr "name" attribute) changes.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-FirstInlineShape($range, $newName) {
    # Re-seating the shape through its own Range before writing Name
    # is required for shapes living in footer stories in this host;
    # it is harmless (and still correct) for header stories too.
    $shape = $range.InlineShapes.Item(1)
    $shape = $shape.Range.InlineShapes.Item(1)
    $shape.Name = $newName
}

# Primary footer -> word/footer2.xml (docPr id="2")
$footerPrimary = $sec.Footers.Item(1)
if ($footerPrimary.Exists -and $footerPrimary.Range.InlineShapes.Count -ge 1) {
    Rename-FirstInlineShape $footerPrimary.Range "image2.png"
}

# First-page footer -> word/footer1.xml (docPr id="3")
$footerFirst = $sec.Footers.Item(2)
if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -ge 1) {
    Rename-FirstInlineShape $footerFirst.Range "image2.png"
}

# First-page header -> word/header1.xml (docPr id="1")
$headerFirst = $sec.Headers.Item(2)
if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -ge 1) {
    Rename-FirstInlineShape $headerFirst.Range "image1.jpg"
}

Write-Output "Renamed inline shapes in footers/header."
